$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text / non-numeric-looking values: safe to assign directly as strings.
$ws.Range("D2").Value = "65.909.09"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "3.201.71"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  +3.65%  "
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.199.12"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("E12").Value = "  +1.58%  "
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  +4.99%  "
$ws.Range("D15").Value = "3.727.96"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D17").Value = "65.994.71"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "3.214.45"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("E21").Value = "  +3.67%  "
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("E23").Value = "  +4.15%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("E29").Value = "  +3.42%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E30").Value = "  +8.69%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("E41").Value = "  +2.45%  "
$ws.Range("E42").Value = "  +6.16%  "
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("D44").Value = "2.951.17"
$ws.Range("E44").Value = "  -4.22%  "
$ws.Range("D45").Value = "0.0₃0644"
$ws.Range("E45").Value = "  +6.15%  "
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  +1.84%  "
$ws.Range("E51").Value = "  +3.92%  "

# Numeric-looking price strings must stay text (they use "." as a thousands
# separator, e.g. "3.201.71"), so force text format before assigning, then
# clear the temporary formatting so no stray style survives in the output.
$numericTextCells = @("D5", "D6", "D9", "D11", "D12", "D13", "D14", "D16", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D46", "D47", "D50", "D51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D5").Value = "599.91"
$ws.Range("D6").Value = "152.91"
$ws.Range("D9").Value = "0.532"
$ws.Range("D11").Value = "6.09"
$ws.Range("D12").Value = "0.511"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("D14").Value = "39.42"
$ws.Range("D16").Value = "7.48"
$ws.Range("D20").Value = "510.91"
$ws.Range("D21").Value = "15.42"
$ws.Range("D22").Value = "0.738"
$ws.Range("D23").Value = "8.15"
$ws.Range("D24").Value = "15.37"
$ws.Range("D25").Value = "84.85"
$ws.Range("D27").Value = "9.29"
$ws.Range("D28").Value = "3.01"
$ws.Range("D29").Value = "2.27"
$ws.Range("D30").Value = "6.86"
$ws.Range("D31").Value = "2.86"
$ws.Range("D32").Value = "28.07"
$ws.Range("D35").Value = "6.56"
$ws.Range("D36").Value = "54.98"
$ws.Range("D37").Value = "0.0905"
$ws.Range("D38").Value = "485.41"
$ws.Range("D39").Value = "0.0420"
$ws.Range("D40").Value = "2.94"
$ws.Range("D41").Value = "8.89"
$ws.Range("D42").Value = "0.303"
$ws.Range("D43").Value = "0.120"
$ws.Range("D46").Value = "2.44"
$ws.Range("D47").Value = "28.50"
$ws.Range("D50").Value = "2.30"
$ws.Range("D51").Value = "2.55"
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).ClearFormats()
}
